$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New BVT row (row 8): "Opacity" feature checklist entry.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Opacity"
$ws.Range("C8").Value = "Update Oparcity"
$ws.Range("D8").Value = "1. Go to formatting pane`n2. Go to Opacity option`n3. Update 'External arcs' to '90'`n4. Update 'Internal arcs' to '25'"
$ws.Range("E8").Value = "1. Opacity of External arcs will be set to '90%'`n2. Opacity of Internal arcs will be set to '25%'"

# Match formatting used by the other BVT rows: thin border all round,
# and wrap text on the Steps/Output columns.
$ws.Range("A8:E8").Borders.LineStyle = 1
$ws.Range("D8:E8").WrapText = $true

$ws.Rows.Item(8).RowHeight = 60

# Move the active selection down to the new row's BVT cell.
$ws.Range("C7").Select() | Out-Null
